# Daily attendance processing - 2025-10-23 21:41:21
# Normalize the "Recorded By" (column G) list order: for entries that were
# recorded starting with "System" (but excluding the admin@admin.com
# entries), reverse the comma-separated order so "System" moves to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "System*") { continue }
    if ($val -like "*admin@admin.com*") { continue }

    $parts = $val -split ",\s*"
    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }
    $newVal = [string]::Join(", ", $reversed)

    $cell.Value2 = $newVal
}
